$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.14"
$ws.Range("D2").Style = "Normal"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "11"
$ws.Range("G2").Style = "Normal"

$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "11"
$ws.Range("G3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.426"
$ws.Range("D4").Style = "Normal"

$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "11"
$ws.Range("G4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05838"
$ws.Range("D5").Style = "Normal"

$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "11"
$ws.Range("G5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.385"
$ws.Range("D6").Style = "Normal"

$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "11"
$ws.Range("G6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.350"
$ws.Range("D7").Style = "Normal"

$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "11"
$ws.Range("G7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8085"
$ws.Range("D8").Style = "Normal"

$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "11"
$ws.Range("G8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9723"
$ws.Range("D9").Style = "Normal"

$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "11"
$ws.Range("G9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1428"
$ws.Range("D10").Style = "Normal"

$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "11"
$ws.Range("G10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07437"
$ws.Range("D11").Style = "Normal"

$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "11"
$ws.Range("G11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03283"
$ws.Range("D12").Style = "Normal"

$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "11"
$ws.Range("G12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03013"
$ws.Range("D13").Style = "Normal"

$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "11"
$ws.Range("G13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.139"
$ws.Range("D14").Style = "Normal"

$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "11"
$ws.Range("G14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09409"
$ws.Range("D15").Style = "Normal"

$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "11"
$ws.Range("G15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001588"
$ws.Range("D16").Style = "Normal"

$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "11"
$ws.Range("G16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04822"
$ws.Range("D17").Style = "Normal"

$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "11"
$ws.Range("G17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005890"
$ws.Range("D18").Style = "Normal"

$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "11"
$ws.Range("G18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006139"
$ws.Range("D19").Style = "Normal"

$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "11"
$ws.Range("G19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004111"
$ws.Range("D20").Style = "Normal"

$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "11"
$ws.Range("G20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009977"
$ws.Range("D21").Style = "Normal"

$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "11"
$ws.Range("G21").Style = "Normal"

$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "11"
$ws.Range("G22").Style = "Normal"

$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "11"
$ws.Range("G23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.223"
$ws.Range("D24").Style = "Normal"

$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "11"
$ws.Range("G24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3209"
$ws.Range("D25").Style = "Normal"

$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "11"
$ws.Range("G25").Style = "Normal"

$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "11"
$ws.Range("G26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003400"
$ws.Range("D27").Style = "Normal"

$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "11"
$ws.Range("G27").Style = "Normal"

$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "11"
$ws.Range("G28").Style = "Normal"

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "11"
$ws.Range("G29").Style = "Normal"

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "11"
$ws.Range("G30").Style = "Normal"

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "11"
$ws.Range("G31").Style = "Normal"

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "11"
$ws.Range("G32").Style = "Normal"

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "11"
$ws.Range("G33").Style = "Normal"

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "11"
$ws.Range("G34").Style = "Normal"

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "11"
$ws.Range("G35").Style = "Normal"

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "11"
$ws.Range("G36").Style = "Normal"

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "11"
$ws.Range("G37").Style = "Normal"

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "11"
$ws.Range("G38").Style = "Normal"

$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "11"
$ws.Range("G39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03880"
$ws.Range("D40").Style = "Normal"

$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "11"
$ws.Range("G40").Style = "Normal"

$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "11"
$ws.Range("G41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1078"
$ws.Range("D42").Style = "Normal"

$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "11"
$ws.Range("G42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002700"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "42CEJICEJIBestin24h"

$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "11"
$ws.Range("G43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.006389"
$ws.Range("D44").Style = "Normal"

$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "11"
$ws.Range("G44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005614"
$ws.Range("D45").Style = "Normal"

$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "11"
$ws.Range("G45").Style = "Normal"

$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "11"
$ws.Range("G46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4200"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "11"
$ws.Range("G47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1456"
$ws.Range("D48").Style = "Normal"

$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "11"
$ws.Range("G48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("D49").Style = "Normal"

$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "11"
$ws.Range("G49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01010"
$ws.Range("D50").Style = "Normal"

$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "11"
$ws.Range("G50").Style = "Normal"

$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "11"
$ws.Range("G51").Style = "Normal"
